$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D, shifting existing D:K data to E:L
$ws.Columns("D").Insert()

# Copy number formatting/styles from column E (the old D data) into new column D.
# Done in segments that skip the unused separator rows (36 and 78) so that no
# stray row/cell entries get created where none existed before.
$ws.Range("E5:E35").Copy()
$ws.Range("D5:D35").PasteSpecial(-4122)
$ws.Range("E37:E77").Copy()
$ws.Range("D37:D77").PasteSpecial(-4122)
$ws.Range("E79:E102").Copy()
$ws.Range("D79:D102").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Populate new column D with the latest fiscal-period figures
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(8, 4).Value = 450400
$ws.Cells.Item(9, 4).Value = 170500
$ws.Cells.Item(10, 4).Value = 279900
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(14, 4).Value = 700
$ws.Cells.Item(15, 4).Value = 164200
$ws.Cells.Item(17, 4).Value = 370000
$ws.Cells.Item(18, 4).Value = 80400
$ws.Cells.Item(20, 4).Value = 42200
$ws.Cells.Item(21, 4).Value = 286800
$ws.Cells.Item(22, 4).Value = 52800
$ws.Cells.Item(23, 4).Value = 69800
$ws.Cells.Item(24, 4).Value = 0
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(26, 4).Value = 69800
$ws.Cells.Item(27, 4).Value = 67500
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(29, 4).Value = "NA"
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(32, 4).Value = -42200
$ws.Cells.Item(33, 4).Value = 67500
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(35, 4).Value = 67500
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(41, 4).Value = 8400
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(43, 4).Value = 79200
$ws.Cells.Item(44, 4).Value = 0
$ws.Cells.Item(45, 4).Value = 66200
$ws.Cells.Item(46, 4).Value = 0
$ws.Cells.Item(47, 4).Value = 8500
$ws.Cells.Item(48, 4).Value = 2958900
$ws.Cells.Item(49, 4).Value = 22000
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(52, 4).Value = 9300
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(54, 4).Value = 3191200
$ws.Cells.Item(57, 4).Value = 80400
$ws.Cells.Item(58, 4).Value = 0
$ws.Cells.Item(59, 4).Value = 0
$ws.Cells.Item(60, 4).Value = 0
$ws.Cells.Item(61, 4).Value = 1346000
$ws.Cells.Item(62, 4).Value = 0
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(66, 4).Value = 1474600
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(72, 4).Value = 0
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(76, 4).Value = 1716600
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(81, 4).Value = 67500
$ws.Cells.Item(83, 4).Value = 164200
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(89, 4).Value = 208400
$ws.Cells.Item(91, 4).Value = -175100
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(94, 4).Value = -105000
$ws.Cells.Item(96, 4).Value = -150300
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(100, 4).Value = -101200
$ws.Cells.Item(101, 4).Value = 0
$ws.Cells.Item(102, 4).Value = 2200

# Apply restated figures to historical columns that were revised
$ws.Cells.Item(8, 5).Value = 424700
$ws.Cells.Item(8, 6).Value = 412000
$ws.Cells.Item(9, 5).Value = 157300
$ws.Cells.Item(10, 5).Value = 267500
$ws.Cells.Item(10, 6).Value = 265400
$ws.Cells.Item(14, 6).Value = 4600
$ws.Cells.Item(15, 6).Value = 255400
$ws.Cells.Item(17, 6).Value = 310100
$ws.Cells.Item(18, 5).Value = 39300
$ws.Cells.Item(18, 6).Value = 101800
$ws.Cells.Item(20, 5).Value = 40200
$ws.Cells.Item(20, 6).Value = 41100
$ws.Cells.Item(21, 6).Value = 270600
$ws.Cells.Item(23, 6).Value = 85600
$ws.Cells.Item(26, 6).Value = 85600
$ws.Cells.Item(27, 6).Value = 85600
$ws.Cells.Item(32, 5).Value = -40200
$ws.Cells.Item(32, 6).Value = -41100
$ws.Cells.Item(33, 6).Value = 85400
$ws.Cells.Item(35, 6).Value = 85400
$ws.Cells.Item(81, 6).Value = 85400
